# ODD materials part6_summary.pptx update:
#  1. Slide 3 ("Learnings"): rewrite the "how to write readable and
#     maintainable code" bullet into a single run with the new wording,
#     and remove the stray red "Oval 1" decoration shape.
#  2. Swap the positions of the "References" slide and the "Buon appetito"
#     (pizza) slide so the pizza/questions slide now comes right after the
#     learnings slide, and References follows it.

$p = $ppt.ActivePresentation

# --- 1a. Update the bullet text on slide 3 -------------------------------
$s3 = $p.Slides.Item(3)
$bodyShape = $s3.Shapes.Item("Text Placeholder 2")
$tr = $bodyShape.TextFrame.TextRange
$para = $tr.Paragraphs(4)

# Force a full run replacement (rather than a partial, prefix-preserving
# edit) so the paragraph collapses back down to a single run.
$para.Text = "x"
$para.Text = "how to write better code by detecting code smells and how to avoid them"

# --- 1b. Remove the red oval decoration shape ----------------------------
$s3.Shapes.Item("Oval 1").Delete()

# --- 2. Reorder slides: move the pizza/"Buon appetito" slide (currently
#        slide 5) to slide position 4, pushing the References slide to 5.
$p.Slides.Item(5).MoveTo(4)
